# Persist score details onto the "Book1" sheet.
#
# Summary of the change:
#  - Fix a typo in the Work Location value for the second employee row
#    ("hyderabad" -> "Hyderabad").
#  - Populate the per-module SCORE columns (L:AD) on row 3, which had been
#    left blank.
#  - The old, now-unused placeholder "overall score" figures in column AL
#    (rows 3-4) are removed now that the detailed per-module scores are
#    persisted instead.
#  - Leave the selection where the user last left it (AH3) when the file
#    was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "hyderabad" -> "Hyderabad" capitalization typo -------------
$ws.Range("J4").Value = "Hyderabad"

# --- Fill in the module score details for row 3 (L3:AD3) ----------------
$ws.Range("L3").Value = 70
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 60
$ws.Range("O3").Value = 90
$ws.Range("P3").Value = 70
$ws.Range("Q3").Value = 70
$ws.Range("R3").Value = 70
$ws.Range("S3").Value = 70
$ws.Range("T3").Value = 60
$ws.Range("U3").Value = 60
$ws.Range("V3").Value = 70
$ws.Range("W3").Value = 65
$ws.Range("X3").Value = 55
$ws.Range("Y3").Value = 55
$ws.Range("Z3").Value = 50
$ws.Range("AA3").Value = 50
$ws.Range("AB3").Value = 90
$ws.Range("AC3").Value = 70
$ws.Range("AD3").Value = 70

# --- Drop the now-superseded overall-score figures in column AL ---------
$ws.Range("AL3").ClearContents()
$ws.Range("AL4").ClearContents()

# --- Restore the last active selection -----------------------------------
$ws.Range("AH3").Select()
